# Daily attendance processing - 2025-12-14 06:36:38
# Normalize the "Recorded By" column (G) so that the current user's
# email (dnasr281@gmail.com) is always listed first in the
# comma-separated list of recorders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetUser = "dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value) { continue }

    $text = [string]$value
    if ($text -notlike "*$targetUser*") { continue }

    $parts = $text -split ", "
    if ($parts.Count -lt 2) { continue }

    # Only reorder when the target user is not already first.
    if ($parts[0] -eq $targetUser) { continue }

    $idx = -1
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($parts[$i] -eq $targetUser) {
            $idx = $i
            break
        }
    }
    if ($idx -lt 0) { continue }

    $rest = @()
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($i -ne $idx) { $rest += $parts[$i] }
    }

    $newParts = @($targetUser) + $rest
    $newText = [string]::Join(", ", $newParts)

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
